# Handles float input without breaking stuff
# Updates the grading/marksheet summary rows and fills in the student-answer
# column with actual results, dropping the unused extra answer sections.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# --- Summary block (rows 10-12) -------------------------------------------
# Row 10: "No." counts
$ws.Range("B10").Value = 14
$ws.Range("C10").Value = 6
$ws.Range("D10").Value = 8
$ws.Range("E10").Value = 28

# Row 11: "Marking" scheme (penalty is numeric, not text, now)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12: "Total"
$ws.Range("B12").Value = 56
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "50/112"

# A10/A11/A12 pick up the bold header style (same as A9/A15) without
# altering their existing label text.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Drop the third Student Ans / Correct Ans block (columns G:H) ---------
$ws.Range("G15:H40").Clear()

# --- Drop the second block (columns D:E) for all but the first three rows -
$ws.Range("D19:E40").Clear()

# --- Fill in student answers (column A) for rows 16-40 --------------------
# Helper style source cells: B10 carries the "correct" (green) style,
# C10 carries the "incorrect" (red) style already used in the sheet.

function Set-Answer($cell, $text, $correct) {
    $ws.Range($cell).Value = $text
    if ($correct) {
        $ws.Range("B10").Copy()
    } else {
        $ws.Range("C10").Copy()
    }
    $ws.Range($cell).PasteSpecial(-4122)
    $excel.CutCopyMode = 0
}

Set-Answer "A16" "Option A" $true
Set-Answer "A17" "Option D" $true
Set-Answer "A18" "Option B" $true
Set-Answer "A19" "Option C" $true
Set-Answer "A21" "Option D" $false
Set-Answer "A23" "Option C" $false
Set-Answer "A26" "Option C" $true
Set-Answer "A27" "Option A" $true
Set-Answer "A28" "Option D" $true
Set-Answer "A29" "Option D" $true
Set-Answer "A31" "Option D" $true
Set-Answer "A33" "Option D" $true
Set-Answer "A35" "Option D" $true
Set-Answer "A37" "Option B" $false
Set-Answer "A38" "Option B" $false
Set-Answer "A39" "Option D" $true
Set-Answer "A40" "Option C" $false

# --- Update the remaining second-block (D:E) rows --------------------------
Set-Answer "D16" "Option A" $true
Set-Answer "D17" "Option C" $true
Set-Answer "D18" "Option B" $false
